
$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate the existing "总计" sheet, placing the copy right after it ---
# The ORIGINAL sheet object keeps its sheetId/rId and becomes the new "2022-Q1" sheet;
# the COPY becomes the new "总计" sheet (with the 2022-Q1 summary row added).
$totalOrig = $wb.Worksheets.Item("总计")
$totalOrig.Copy($null, $totalOrig)
$totalCopy = $wb.ActiveSheet

$totalOrig.Name = "2022-Q1"
$totalCopy.Name = "总计"

# =========================================================================
# Step 2: rebuild the "2022-Q1" sheet ($totalOrig) as a fund-holdings table
# =========================================================================
$q1 = $totalOrig

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("G1").Value = "持有市值(亿元)"

# New header cells E1, F1, H1 need the same style as the existing header cells;
# copy the formatting from D1 (style index already used by B1:D1) onto them.
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("H1").Value = "仓位排名"

# Fund holdings data rows (A2:H23)
$data = @(
    @('003293','易方达科瑞灵活配置混合','34.67','78.17','4.12','1.4284',2),
    @('006533','易方达科融混合','32.14','89.61','2.41','0.7746',7),
    @('110012','易方达科汇灵活配置混合','15.73','75.64','4.00','0.6292',3),
    @('110002','易方达策略成长混合','12.15','88.76','4.75','0.5771',3),
    @('112002','易方达策略成长二号混合','10.72','87.99','4.32','0.4631',3),
    @('519909','华安安顺灵活配置混合','10.93','87.45','3.02','0.3301',9),
    @('001070','建信信息产业股票','10.67','81.80','2.84','0.3030',10),
    @('011649','易方达逆向投资混合A','7.49','85.02','3.63','0.2719',5),
    @('010389','易方达科益混合A','7.10','92.94','3.70','0.2627',9),
    @('005310','广发电子信息传媒产业精选股票A','3.99','90.16','4.13','0.1648',6),
    @('000308','建信创新中国混合','3.11','84.50','3.10','0.0964',9),
    @('011650','易方达逆向投资混合C','1.96','85.02','3.63','0.0711',5),
    @('159804','国寿安保国证创业板中盘精选88ETF','2.10','98.79','2.25','0.0472',4),
    @('001534','华宝万物互联灵活配置混合','1.06','92.81','3.54','0.0375',8),
    @('010236','广发电子信息传媒产业精选股票C','0.81','90.16','4.13','0.0335',6),
    @('160812','长盛同益成长回报灵活配置混合(LOF)','1.49','81.36','2.21','0.0329',6),
    @('002152','华宝核心优势灵活配置混合','0.45','90.91','3.44','0.0155',8),
    @('002789','长盛同享灵活配置混合A','0.49','79.17','2.38','0.0117',7),
    @('010390','易方达科益混合C','0.29','92.94','3.70','0.0107',9),
    @('000892','九泰天宝灵活配置混合A','0.07','90.81','4.51','0.0032',10),
    @('002790','长盛同享灵活配置混合C','0.02','79.17','2.38','0.0005',7),
    @('002028','九泰天宝灵活配置混合C','0.00','90.81','4.51','0',10)
)

# Columns B:G hold text (to preserve formatting such as trailing zeros); force
# text storage by setting NumberFormat to Text before assignment, then restore
# the default "Normal" style so no residual number formatting is left behind.
$q1.Range("B2:G23").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $q1.Cells.Item($r, 1).Value = ($r - 2)
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r++
}

$q1.Range("B2:G23").Style = "Normal"

# The last row's 持有市值 (G23) is a genuine number (0), not text - restore it.
$q1.Range("G23").NumberFormat = "General"
$q1.Range("G23").Value = 0
$q1.Range("G23").Style = "Normal"

# A2:A23 index column + H2:H23 rank column match the style already used by A2
# (copied along with the sheet) - make sure every row in the extended range
# carries that same index style.
$q1.Range("A2").Copy()
$q1.Range("A2:A23").PasteSpecial(-4122)

# =========================================================================
# Step 3: rebuild the "总计" sheet ($totalCopy) with the new 2022-Q1 summary row
# =========================================================================
$tot = $totalCopy

# Shift the five existing summary rows down by one (row 6 -> 7, ... row 2 -> 3)
$tot.Cells.Item(7, 1).Value = 5
$tot.Cells.Item(7, 2).Value = "2020-Q4"
$tot.Cells.Item(7, 3).Value = 30
$tot.Cells.Item(7, 4).Value = 13.47

$tot.Cells.Item(6, 1).Value = 4
$tot.Cells.Item(6, 2).Value = "2021-Q1"
$tot.Cells.Item(6, 3).Value = 17
$tot.Cells.Item(6, 4).Value = 9.16

$tot.Cells.Item(5, 1).Value = 3
$tot.Cells.Item(5, 2).Value = "2021-Q2"
$tot.Cells.Item(5, 3).Value = 23
$tot.Cells.Item(5, 4).Value = 9.24

$tot.Cells.Item(4, 1).Value = 2
$tot.Cells.Item(4, 2).Value = "2021-Q3"
$tot.Cells.Item(4, 3).Value = 12
$tot.Cells.Item(4, 4).Value = 6.1

$tot.Cells.Item(3, 1).Value = 1
$tot.Cells.Item(3, 2).Value = "2021-Q4"
$tot.Cells.Item(3, 3).Value = 16
$tot.Cells.Item(3, 4).Value = 2.67

# New top row: 2022-Q1 summary
$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 22
$tot.Cells.Item(2, 4).Value = 5.57

# Row 7 (old row 6, "2020-Q4") is a brand-new row in this sheet's used range,
# so its index cell needs the same style as the other index cells (A2:A6).
$tot.Range("A2").Copy()
$tot.Range("A7").PasteSpecial(-4122)

# Restore the originally-active sheet/tab (the Copy() calls above shifted the
# active tab onto the duplicated sheet).
$wb.Worksheets.Item("2020-Q4").Activate()
